# developer_week_01.xlsx - "add .dev and estructura_tablas.xls"
#
# The functional change on the "Hoja1" worksheet is a new column inserted
# before column E. That pushes the old E:I block (the pronunciation-audio
# marker columns plus the "repaso" word list) one column to the right, to
# F:J, while A:D (word / translation / pronunciation / audio columns) stay
# put. The new column inherits the column width/format of its left
# neighbour (column D), which is exactly what Excel's normal "insert
# column" does.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Insert the new column at E - shifts E:I -> F:J.
$ws.Columns("E:E").Insert()

# Excel's column-insert defaults to the formatting of the column to the
# left; make that explicit so the new column's width matches column D's,
# exactly like the committed workbook shows for cols D:E (both sized the
# same as the original column D).
$ws.Columns("E:E").ColumnWidth = 17.5

# Restore the active selection recorded in the saved view: D2:D42 selected
# with D2 as the active cell, and the frozen-header pane scrolled down so
# row 31 is the first visible row below the freeze line.
$ws.Range("A31").Select()
$ws.Range("D2:D42").Select()
